$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5.968389987945557
$ws.Range("B1").Value = 3.222516775131226
$ws.Range("C1").Value = 1.947881698608398
$ws.Range("D1").Value = 1.623064398765564
$ws.Range("E1").Value = 1.533383727073669
